# Opentrons 02-ont-gDNA template: switch the "Used here" math from a
# 7 Mb plasmid genome (0.023 fmoles) to a 30 kb gDNA fragment (5.4 fmoles)
# so reaction volumes flex with the new target size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Explanatory note next to the constants block (I3)
$ws.Range("I3").Value = "Used here - 100 ng @ 30 kb DNA genome = 5.4 fmoles"

# "fmoles per rxn" constant (J5): 0.04 -> 5.4
$ws.Range("J5").Value = 5.4

# "plasmid size" column (E) for every sample row (9-25): 5,000,000 bp -> 30,000 bp
for ($row = 9; $row -le 25; $row++) {
    $ws.Cells.Item($row, 5).Value = 30000
}

# The "rapid barcode plate" label shifts one column left, from R8 to Q8,
# now that column R is no longer part of the used range.
$ws.Range("Q8").Value = $ws.Range("R8").Value2
$ws.Range("R8").ClearContents()

# Reflect the new focal cell in the saved view/selection.
$ws.Range("Q8").Select() | Out-Null
